$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "rules"
$ws.Range("D1").Value = "adaptive_filter"

# Update D column (RLS_option -> adaptive_filter) values to text "wRLS"
# and refresh RMSE/NDEI/MAE values with updated precision
$ws.Range("D2").Value = "wRLS"
$ws.Range("E2").Value = 6.025546986248139
$ws.Range("F2").Value = 1.005738801678266
$ws.Range("G2").Value = 4.798328146675348

$ws.Range("D3").Value = "wRLS"
$ws.Range("E3").Value = 6.150385157906561
$ws.Range("F3").Value = 1.026575846589541
$ws.Range("G3").Value = 4.640152593538985

$ws.Range("D4").Value = "wRLS"
$ws.Range("E4").Value = 6.005620046969991
$ws.Range("F4").Value = 1.002412747450085
$ws.Range("G4").Value = 4.340876560880753

$ws.Range("D5").Value = "wRLS"
$ws.Range("E5").Value = 6.074465594736091
$ws.Range("F5").Value = 1.013903926403495
$ws.Range("G5").Value = 4.394404973155996

$ws.Range("D6").Value = "wRLS"
$ws.Range("E6").Value = 5.986518648599366
$ws.Range("F6").Value = 0.9992244862762033
$ws.Range("G6").Value = 4.367646362440475

$ws.Range("D7").Value = "wRLS"
$ws.Range("E7").Value = 5.9552208353104
$ws.Range("F7").Value = 0.9940004916240607
$ws.Range("G7").Value = 4.694499556510445

$ws.Range("D8").Value = "wRLS"
$ws.Range("E8").Value = 5.932163053369541
$ws.Range("F8").Value = 0.9901518607808386
$ws.Range("G8").Value = 4.641095578916298
